# Aggiunto telegramma per il fine curva
# - Aggiornato excel telegrammi
# - Aggiunto telegramma valore #8 per il reset delle variabili cinematiche

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valori")

# Row 10 was an empty placeholder row in the "Valori" table; fill it in
# with the new "Fine curva/reset" telegram (value #8).
$ws.Range("A10").Value = "Fine curva/reset"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = "0x08"

# Move the selection to D11, matching the saved selection state.
$ws.Range("D11").Select()
